# Add two new "Title Only" slides after the existing slide 1, each
# carrying just a title placeholder with the note's title text.
#
# ppLayoutTitleOnly = 11 (maps to the "Title Only" slide layout, the
# layout that only exposes a title placeholder - matches the new
# slide2.xml / slide3.xml shape trees in the target deck).

$p = $ppt.ActivePresentation

$s2 = $p.Slides.Add(2, 11)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "ADBL Banking Note"

$s3 = $p.Slides.Add(3, 11)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "TSC Note"
